# Insert a new daily-push row at row 597 (date 2026/01/09, 金, time 8, rank 201),
# shifting all subsequent rows (old 597..638) down by one (to 598..639).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing row 597 (and everything below it) down by one row.
$ws.Rows.Item(597).Insert()

# Fill in the newly inserted row. The leading apostrophe forces the date-like
# string to be kept as text instead of being auto-converted to a date serial,
# matching the existing "日付" column's text storage; resetting the style back
# to Normal afterwards avoids leaving a stray quote-prefixed style on the cell.
$ws.Range("A597").Value = "'2026/01/09"
$ws.Range("A597").Style = "Normal"
$ws.Range("B597").Value = "金"
$ws.Range("C597").Value = 8
$ws.Range("D597").Value = 201
